# This script reproduces a manual editing session against Data.xlsx:
#  - a couple of login/e-mail values were retyped on a few sheets,
#  - column A on the "Login" sheet was widened,
#  - the user ended up leaving the selection on the "LoginBo" sheet
#    (which becomes the active tab), with a few other sheets' selections
#    left where the user last clicked.

$wb = $excel.ActiveWorkbook

# --- Login sheet: update a couple of e-mail addresses ---------------------
$wsLogin = $wb.Worksheets.Item("Login")
$wsLogin.Range("A3").Value = "rituparna+hotel@needleinnovision.com"
$wsLogin.Range("A1").Value = "rakesh@gmail.com"

# widen column A to fit the new text (closest the host's pixel-snapped
# column-width grid gets to the original 44.85546875 stored width)
$wsLogin.Columns.Item(1).ColumnWidth = 44

# leave the selection where the user last clicked on this sheet
[void]$wsLogin.Range("A6").Select()

# --- LoginFHE sheet: update e-mail address ---------------------------------
$wsLoginFHE = $wb.Worksheets.Item("LoginFHE")
$wsLoginFHE.Range("A1").Value = "ramesh@gmail.com"
[void]$wsLoginFHE.Range("A7").Select()

# --- LoginBo sheet: update e-mail address, end up as the active sheet -----
$wsLoginBo = $wb.Worksheets.Item("LoginBo")
$wsLoginBo.Range("A1").Value = "mandeepm+bomumb@needleinnovision.com"
[void]$wsLoginBo.Activate()
[void]$wsLoginBo.Range("A9").Select()
